$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '62.124.35'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -0.21%  '

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.429.50'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +0.03%  '

# Row 4
$ws.Range("E4").Value = '  -0.02%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '410.10'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.79%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '130.50'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -0.03%  '

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.635'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +6.65%  '

# Row 8
$ws.Range("E8").Value = '  -0.07%  '

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.741'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +7.19%  '

# Row 10
$ws.Range("E10").Value = '  +5.47%  '

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '43.01'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +2.51%  '

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.0000226'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +51.45%  '

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '9.23'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +10.32%  '

# Row 14
$ws.Range("E14").Value = '  -0.19%  '

# Row 15
$ws.Range("B15").Value = 'Chainlink'
$ws.Range("C15").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '21.40'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +7.80%  '

# Row 16
$ws.Range("B16").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C16").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '3.974.86'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +0.00%  '

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '3.426.87'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +0.42%  '

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '12.53'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +8.23%  '

# Row 19
$ws.Range("E19").Value = '  +7.97%  '

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '62.092.43'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -0.21%  '

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '458.26'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +47.04%  '

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '91.86'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +9.40%  '

# Row 23
$ws.Range("E23").Value = '  +1.85%  '

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '13.17'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +3.02%  '

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '3.30'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +4.48%  '

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '33.26'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +12.23%  '

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '9.18'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +12.88%  '

# Row 28
$ws.Range("E28").Value = '  +1.03%  '

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '7.69'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -1.62%  '

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '12.11'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +6.75%  '

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '2.70'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -2.03%  '

# Row 32
$ws.Range("E32").Value = '  -0.60%  '

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '43.16'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -3.38%  '

# Row 34
$ws.Range("E34").Value = '  -0.56%  '

# Row 35
$ws.Range("E35").Value = '  -0.08%  '

# Row 36
$ws.Range("E36").Value = '  +3.67%  '

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '54.38'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +4.85%  '

# Row 38
$ws.Range("E38").Value = '  -0.08%  '

# Row 39
$ws.Range("E39").Value = '  +2.02%  '

# Row 40
$ws.Range("E40").Value = '  +8.00%  '

# Row 41
$ws.Range("E41").Value = '  -0.76%  '

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.320'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -0.74%  '

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '142.50'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +0.31%  '

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '4.29'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +9.91%  '

# Row 45
$ws.Range("B45").Value = 'ARBITRUM'
$ws.Range("C45").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '2.01'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +1.67%  '

# Row 46
$ws.Range("B46").Value = 'WEMIXToken'
$ws.Range("C46").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.55'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +14.93%  '

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '16.69'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -0.54%  '

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '22.45'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +5.55%  '

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.143'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +20.49%  '

# Row 50
$ws.Range("B50").Value = 'ThetaToken'
$ws.Range("C50").Value = 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '2.12'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +7.78%  '

# Row 51
$ws.Range("B51").Value = 'RocketPoolETH'
$ws.Range("C51").Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '3.778.20'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -0.14%  '
